$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'Unique record identifier generated by MAP based upon the YEAR, MONTH, INCIDENT and ORI code in the report. Alphanumeric 16-character format (A16).'
$ws.Range("C3").Value = 'The alphanumeric variable describing the name of the law enforcement agency making the report. (A450).'
$ws.Range("C4").Value = '– The one-digit numeric code describing the type of law enforcement agency making the report so that 1=Sheriff, 2=County Police, 3=Municipality, 5=Primary State Law Enforcement usually meaning the State Police, 6=Special Police, 7=Constable, 8=Tribal Police, 9=Regional Police. When using the Comma Separated Values format file (file extension=.csv) the original FBI numbering scheme is replaced with the alphanumeric label. (F1.0).'
$ws.Range("C7").Value = 'Alphanumeric variable describing the original FBI naming and abbreviating scheme for the state of the reporting agency. (A6).'
$ws.Range("C5").Value = 'The one-digit numeric code describing the type of law enforcement agency making the report so that 1=Sheriff, 2=County Police, 3=Municipality, 5=Primary State Law Enforcement usually meaning the State Police, 6=Special Police, 7=Constable, 8=Tribal Police, 9=Regional Police. When using the Comma Separated Values format file (file extension=.csv) the original FBI numbering scheme is replaced with the alphanumeric label. (F1.0).'
$ws.Range("C8").Value = 'Year of homicide (or when victim’s body was recovered.) Numeric four digit. (F4.0).'
$ws.Range("C10").Value = '– A three-digit number describing the case number within the month in which a homicide occurred. This does not necessarily correspond to the actual case number used inhouse by police agencies. It is used to assist in building a unique record number for each case and to differentiate each case reported within the same month. (F3.0)'
$ws.Range("C11").Value = 'An alphanumeric variable defining whether the report was “A” = “Murder or Nonnegligent manslaughter” or “B” = “Manslaughter by Negligence.” (A1).'
$ws.Range("C12").Value = 'MAP-generated indicator whether Offender was identified at time report was made (SOLVED=1) or not identified (SOLVED=0). Numeric single digit format. (F1.0).'
$ws.Range("C13").Value = 'An alphanumeric variable representing whether the victim was “M” = ”Male” or “F” = ”Female” or “U” indicating “Unknown” gender, usually for conditions in which incomplete remains were recovered. (A1).'
$ws.Range("C14").Value = 'A three-digit numeric variable describing the age in years of the victim. To allow for simpler mathematical calculations, MAP has changed the original alphanumeric coding of “NB” for new born and “BB” for infant to a numeric value of zero to indicate the victim had not achieved a full year of life. A value of 99, as in the original numbering scheme, represents all victims 99 or older. A value of 999 represents victims whose age was not reported, usually because the victim was unidentified and the age was unknown. (F3.0).'
$ws.Range("C15").Value = 'An alphanumeric variable representing whether the victim was “A” = “Asian or Pacific Islander” or “B” = “Black” or “I” = “American Indian or Alaskan Native” or “W” = “White” or “U” victim was of “Unknown” race. (A1).'
$ws.Range("C16").Value = 'An alphanumeric variable representing whether the victim was “H” = “Hispanic Origin” or “N” = “Not of Hispanic Origin” or “U” = “Unknown or Not Reported.” It should be noted that many agencies decline reporting the ethnicity of victims and offenders. (A1).'
$ws.Range("C17").Value = '– An alphanumeric variable representing whether the offender was “M”=”Male” or “F”=”Female” or “U” indicating “Unknown” gender, usually in conditions in which the offender had not been identified at the time of the report. (A1).'
$ws.Range("C18").Value = 'A three-digit numeric variable describing the age in years of the offender. When the offender was not identified at the time of the report, age was reported as 999. A value of 99 represents all offenders 99 or older. (F3.0).'
$ws.Range("C19").Value = '– An alphanumeric variable representing whether the offender was “A” = “Asian or Pacific Islander” or “B” = “Black” or “I” = “American Indian or Alaskan Native” or “W” = “White” or “U” = “Unknown” race, usually in conditions in which the offender had not been identified at the time of the report. (A1).'
$ws.Range("C20").Value = '– An alphanumeric variable representing whether the offender was “H” = “Hispanic Origin” or “N” = “Not of Hispanic Origin” or “U” = “Unknown or Not Reported.” It should be noted that many agencies decline reporting the ethnicity of victims and offenders. (A1).'
$ws.Range("C22").Value = '– A two-digit numeric variable representing the weapon used in the crime. (F2.0).'
$ws.Range("C21").Value = 'An alphanumeric variable describing the relationship between the victim and the offender, if any. '
$ws.Range("C23").Value = 'The number of additional victims (not counting the victim included in the current record) included in the Supplementary Homicide Report’s incident record, which can accept up to 10 additional victims in a single incident report. (F3.0) Incidents of mass murder of more than 11 victims would require multiple SHR incident reports. The Murder Accountability Project’s database captures all reported homicide victims as separate cases. Associated victims will have identical ID numbers'
$ws.Range("C24").Value = 'The number of additional offenders (not counting the offender included in the current record) included in the Supplementary Homicide Report’s incident record, which can accept up to 10 additional offenders in a single incident report. (F3.0) Unlike victims, the Murder Accountability Project does not create multiple case reports for additional offenders. Only the first offender listed in the original SHR report is included.'
$ws.Range("C6").Value = 'An eight-digit numeric variable representing the Census Bureau’s Federal Information Processing Standards (FIPS) code for the Metropolitan Statistical Area from which a record was reported. When using the Comma Separated Values format file (file extension=.csv) the original FIPS coding is replaced with the label indicating the metropolitan area. (F8.0).'
$ws.Range("C9").Value = 'The month of homicide occurrence or when the victim’s body was recovered. '
$ws.Range("C25").Value = 'MAP-generated identifier if record provided by FBI (SOURCE=1) or was obtained by MAP under the Freedom of Information Act from an Agency not participating in SHR reporting to the FBI (SOURCE=0). Numeric single digit format. (F1.0).'

# The "Relationship" description is long, so wrap the text and grow the row
# to fit it (matches the Time Series plot fix picking up full descriptions).
$ws.Range("C21").WrapText = $true
$ws.Rows(21).RowHeight = 210

# Restore the author's last-saved cursor position.
[void]$ws.Range("F15").Select()
